# Weekly update: a new week of price data (2023-01-05) is published for
# "Terminal La Palmera de La Serena - Plátano". The sheet keeps data in
# blocks of 3 rows (Pintón / Primera Maduro / Primera Pintón) ordered with
# the newest week at the top of the data range (row 885) followed by all
# prior weeks. Publishing the new week pushes the existing rows down by
# one block (3 rows) and inserts the fresh block at the top.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at the top of the data block (row 885), shifting
# the existing rows 885:947 down to 888:950. This also grows the sheet
# dimension from T947 to T950 and carries the column styles (e.g. the date
# format on column D) down to the new rows, matching Excel's native
# "Insert Copied Cells"/"Insert Rows" behaviour.
$ws.Range("A885:A887").EntireRow.Insert()

$newDate = 44931

$newRows = @(
    @("Pintón",          80,  18000, 900),
    @("Primera Maduro",  120, 20000, 1000),
    @("Primera Pintón",  120, 21000, 1050)
)

for ($i = 0; $i -lt 3; $i++) {
    $r = 885 + $i
    $vol   = $newRows[$i][1]
    $price = $newRows[$i][2]
    $kg    = $newRows[$i][3]

    $ws.Cells.Item($r, 1).Value  = 8
    $ws.Cells.Item($r, 2).Value  = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value  = "Coquimbo"
    $ws.Cells.Item($r, 4).Value  = $newDate
    $ws.Cells.Item($r, 5).Value  = 4
    $ws.Cells.Item($r, 6).Value  = "Fruta"
    $ws.Cells.Item($r, 7).Value  = 100108
    $ws.Cells.Item($r, 8).Value  = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value  = 100108006
    $ws.Cells.Item($r, 10).Value = "Plátano"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 13).Value = $vol
    $ws.Cells.Item($r, 14).Value = $price
    $ws.Cells.Item($r, 15).Value = $price
    $ws.Cells.Item($r, 16).Value = $price
    $ws.Cells.Item($r, 17).Value = "$/caja 20 kilos"
    $ws.Cells.Item($r, 18).Value = "Ecuador"
    $ws.Cells.Item($r, 19).Value = $kg
    $ws.Cells.Item($r, 20).Value = 20
}
